$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain text so numeric-looking strings
# (e.g. "37.199.75", "57.30") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.199.75'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '2.060.20'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '248.95'
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '57.30'
$ws.Range("E8").Value = '  -2.16%  '
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = '0.915'
$ws.Range("E13").Value = '  +13.83%  '
$ws.Range("D14").Value = '2.361.76'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = '5.77'
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("D16").Value = '2.060.40'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = '18.84'
$ws.Range("E17").Value = '  +12.35%  '
$ws.Range("D18").Value = '37.243.15'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = '5.49'
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("D22").Value = '238.25'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '2.49'
$ws.Range("E24").Value = '  +4.76%  '
$ws.Range("E25").Value = '  +4.35%  '
$ws.Range("E26").Value = '  -4.17%  '
$ws.Range("D27").Value = '170.36'
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").Value = '5.18'
$ws.Range("E30").Value = '  +9.58%  '
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '4.69'
$ws.Range("E32").Value = '  +5.33%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0627'
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").Value = '0.0892'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("E39").Value = '  +14.35%  '
$ws.Range("D40").Value = '3.10'
$ws.Range("E40").Value = '  +8.10%  '
$ws.Range("D41").Value = '0.101'
$ws.Range("E41").Value = '  -11.55%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0224'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '17.64'
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = '96.91'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = '1.278.90'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").Value = '2.86'
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").Value = '2.251.17'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").Value = '44.56'
$ws.Range("E51").Value = '  +1.89%  '

# Restore default (General) number format/style for Price column so cells
# match the original unstyled appearance.
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"

